$d = $word.ActiveDocument

# Locate the sentence that needs the en dash inserted inside the word
# "imprime" ("i" + EN DASH + "mprime"), matching the target edit exactly.
$search = "imprime uma linha em branco"
$rng = $d.Content
$found = $rng.Find.Execute($search, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $beforeRevisionCount = $d.Revisions.Count

    # Enable track-changes so the inserted character is recorded as its own
    # revision. Accepting that single revision afterwards forces Word to keep
    # the new text in its own run (instead of silently re-merging it back
    # into the surrounding run the way a plain InsertAfter would).
    $wasTracking = $d.TrackRevisions
    $d.TrackRevisions = $true

    # $rng.Start now points right before the "i" of "imprime"; collapse a
    # range to just after that "i" and insert the en dash (U+2013) there.
    $insertPoint = $d.Range($rng.Start + 1, $rng.Start + 1)
    $insertPoint.InsertAfter([string][char]0x2013)

    $d.TrackRevisions = $wasTracking

    # Only accept the revision(s) that our edit just introduced, leaving any
    # pre-existing tracked changes in the document untouched.
    for ($i = $d.Revisions.Count; $i -ge ($beforeRevisionCount + 1); $i--) {
        $d.Revisions.Item($i).Accept()
    }
}
